# costs_inpt.xlsx — inpatient and outpatient costs for state S1 (row 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: inpatient cost for S1 = ROUND(1.2*(1375.84+121.61+93.2),0) -> 1909
$ws.Range("B2").Formula = "=ROUND(1.2*(1375.84+121.61+93.2),0)"

# C2: outpatient cost for S1, converted/deflated like rows 3 & 4 -> 469
$ws.Range("C2").Formula = "=ROUND(B2*((2927.15-1027.33)/3.92)/1971.79,0)"

# D2: reference for the new S1 cost figures
$ws.Range("D2").Value = "Graham2018BIA"

# Row 2 grows to the same wrapped-text height as rows 3 & 4 now that D2 has text
$ws.Rows.Item(2).RowHeight = 17

# Scroll the sheet so row 2 is visible at the top, leaving the C3 selection untouched
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
